$d = $word.ActiveDocument
$apos = [char]0x2019

# First change: expand/replace the "Further requirement is very clear although output label..." sentence
$d.Content.Find.Execute(
    "Further requirement is very clear although output label wasn" + $apos + "t given so its fall under Semi Supervised Learning. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Further requirement and output label is clear, thus its falls under the Supervised Learning. ",
    2)

# Second change: " Semi supervised Learning " -> " Supervised Learning "
$d.Content.Find.Execute(
    "Semi supervised Learning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Supervised Learning",
    2)

# Move the _GoBack bookmark: it used to sit at the end of the "Might Resign"
# paragraph; it now marks the last-edited spot, right before "Supervised
# Learning." in the sentence we just rewrote.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$r = $d.Content
$r.Find.Execute("the Supervised Learning.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.Start + 4
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
